$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected for "Feria Lagunitas de Puerto
# Montt - Albahaca". It belongs right above the existing row 124, so shift
# every row from 124 down by inserting a fresh row (xlShiftDown = -4121),
# then populate that new row 124 with the new record's data.
$ws.Rows(124).Insert(-4121)

$ws.Range("A124").Value = 4
$ws.Range("B124").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C124").Value = "Los Lagos"
$ws.Range("D124").Value = 44992
$ws.Range("E124").Value = 10
$ws.Range("F124").Value = 100112052
$ws.Range("G124").Value = "Albahaca"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 80
$ws.Range("K124").Value = 6500
$ws.Range("L124").Value = 6500
$ws.Range("M124").Value = 6500
$ws.Range("N124").Value = "$/docena de matas"
$ws.Range("O124").Value = "Región Metropolitana"
$ws.Range("P124").Value = 1083
$ws.Range("Q124").Value = 6
$ws.Range("R124").Value = "Hortaliza"
